$d = $word.ActiveDocument

# Avoid Word's "smart quotes" autocorrect mangling straight quotes/apostrophes
# when we programmatically set text.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title
#    paragraph at the top of the document.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.ParagraphFormat.Style = "Normal"

$metaStart = $metaPara.Range.Start
$metaLabel = "Meta description"
$metaFull = $metaLabel + ": Read our unbiased review of Christmas Gold Digger, the festive slot game. Play for free and discover pros, cons, RTP, and special features."

$metaPara.Range.Text = $metaFull

$metaBoldRange = $d.Range($metaStart, $metaStart + $metaLabel.Length)
$metaBoldRange.Font.Bold = 1

# ---------------------------------------------------------------------------
# 2) Near the end of the document, remove the duplicated bold title
#    paragraph ("Play Christmas Gold Digger for Free - Review 2021") that
#    now lives just before the closing italic paragraph.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$dupTitlePara.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new DALL-E
#    image-prompt text, keeping its existing italic run formatting intact.
# ---------------------------------------------------------------------------
$count2 = $d.Paragraphs.Count
$closingPara = $d.Paragraphs.Item($count2)
$closingRange = $closingPara.Range
$closingTextRange = $d.Range($closingRange.Start, $closingRange.End - 1)

$newClosingText = 'Dear DALLE, I need a feature image for the online slot game "Christmas Gold Digger". The image should be in a cartoon style and feature a happy Maya warrior with glasses. The main focus of the image should be on the Maya warrior, with elements of Christmas and gold mining included in the background. The image should look festive and fun, with bright colors and playful details that will appeal to players. Please include the game''s title "Christmas Gold Digger" in the image, as well as any other elements that you think will help players identify and enjoy this game. Thank you!'

$closingTextRange.Text = $newClosingText
